$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# New data for rows 2-46 (Timestamp shifted +13 days, new B/C values, plus 3 extra rows)
$dates = @(
    45944,
    45944.01041666666,
    45944.02083333334,
    45944.03125,
    45944.04166666666,
    45944.05208333334,
    45944.0625,
    45944.07291666666,
    45944.08333333334,
    45944.09375,
    45944.10416666666,
    45944.11458333334,
    45944.125,
    45944.13541666666,
    45944.14583333334,
    45944.15625,
    45944.16666666666,
    45944.17708333334,
    45944.1875,
    45944.19791666666,
    45944.20833333334,
    45944.21875,
    45944.22916666666,
    45944.23958333334,
    45944.25,
    45944.26041666666,
    45944.27083333334,
    45944.28125,
    45944.29166666666,
    45944.30208333334,
    45944.3125,
    45944.32291666666,
    45944.33333333334,
    45944.34375,
    45944.35416666666,
    45944.36458333334,
    45944.375,
    45944.38541666666,
    45944.39583333334,
    45944.40625,
    45944.41666666666,
    45944.42708333334,
    45944.4375,
    45944.44791666666,
    45944.45833333334
)

$colB = @(
    1.613, 9.273, 0.875, 2.297, 3.174, 6.307, 5.57, 0.628, 5.235, 0,
    1.742, 4.598, 1.572, 0, 0, 0, 0, 0.195, 0.076, 0,
    0, 0.038, 0.894, 9.673999999999999, 4.336, 5.706, 6.249, 12.104, 1.099, 0,
    0, 0, 1.203, 0.008, 0, 0, 0.122, 0.003, 0, 0,
    0, 0, 0, 0, 0
)

$colC = @(
    5.582, 0.102, 5.934, 3.881, 0.599, 0.002, 0.459, 4.611, 1.354, 13.796,
    5.616, 0.096, 4.455, 14.954, 12.306, 14.513, 1.754, 1.112, 1.607, 2.564,
    22.897, 8.913, 5.072, 0, 0.005, 0, 0, 0, 3.9, 9.195,
    17.105, 32.792, 2.69, 7.036, 24.745, 8.907, 9.119, 8.893000000000001, 10.166, 10.995,
    3.164, 1.991, 5.587, 2.668, 5.505
)

for ($i = 0; $i -lt $dates.Length; $i++) {
    $row = $i + 2
    $ws.Cells.Item($row, 1).Value = $dates[$i]
    $ws.Cells.Item($row, 2).Value = $colB[$i]
    $ws.Cells.Item($row, 3).Value = $colC[$i]
}

# Apply the date style (style index 2, numFmtId 164) to the new A column cells (rows 44-46)
$srcStyleRange = $ws.Range("A43")
for ($row = 44; $row -le 46; $row++) {
    $ws.Cells.Item($row, 1).NumberFormat = $srcStyleRange.NumberFormat
}
